$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-06 05:17:43'
$ws.Range("O2").Value = '-1.6 °C'
$ws.Range("E3").Value = '2026-02-06 05:17:46'
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '74%'
$ws.Range("I3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("O3").Value = '-2.4 °C'
$ws.Range("E4").Value = '2026-02-06 05:17:48'
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = '59%'
$ws.Range("I4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("J4").Value = '993.1 hPa'
$ws.Range("N4").Value = '9.1 °C 4:59 TU'
$ws.Range("O4").Value = '12.7 °C'
$ws.Range("E5").Value = '2026-02-06 05:17:51'
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '76%'
$ws.Range("I5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("J5").Value = '993.6 hPa'
$ws.Range("N5").Value = '6.3 °C 4:38 TU'
$ws.Range("O5").Value = '8.1 °C'
$ws.Range("E6").Value = '2026-02-06 05:17:53'
$ws.Range("J6").Value = '994.6 hPa'
$ws.Range("K6").Value = '-0.1 MJ/m2'
$ws.Range("O6").Value = '14.4 °C'
$ws.Range("E7").Value = '2026-02-06 05:17:56'
$ws.Range("J7").Value = '994.5 hPa'
$ws.Range("N7").Value = '9.5 °C 4:30 TU'
$ws.Range("O7").Value = '10.0 °C'
$ws.Range("E8").Value = '2026-02-06 05:17:58'
$ws.Range("K8").Value = '-0.1 MJ/m2'
$ws.Range("N8").Value = '4.1 °C 4:59 TU'
$ws.Range("O8").Value = '5.9 °C'
$ws.Range("E9").Value = '2026-02-06 05:18:00'
$ws.Range("N9").Value = '0.6 °C 4:53 TU'
$ws.Range("O9").Value = '2.2 °C'
$ws.Range("E10").Value = '2026-02-06 05:18:03'
$ws.Range("N10").Value = '3.2 °C 4:59 TU'
$ws.Range("O10").Value = '5.1 °C'
$ws.Range("E11").Value = '2026-02-06 05:18:05'
$ws.Range("J11").Value = '995.3 hPa'
$ws.Range("N11").Value = '2.4 °C 4:52 TU'
$ws.Range("O11").Value = '4.5 °C'
$ws.Range("E12").Value = '2026-02-06 05:18:08'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '60%'
$ws.Range("I12").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("N12").Value = '9.1 °C 4:55 TU'
$ws.Range("O12").Value = '12.4 °C'
$ws.Range("E13").Value = '2026-02-06 05:18:10'
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '89%'
$ws.Range("I13").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("N13").Value = '4.0 °C 4:59 TU'
$ws.Range("O13").Value = '6.5 °C'
$ws.Range("E14").Value = '2026-02-06 05:18:13'
$ws.Range("O14").Value = '-3.7 °C'
$ws.Range("E15").Value = '2026-02-06 05:18:15'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '84%'
$ws.Range("I15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("J15").Value = '993.7 hPa'
$ws.Range("N15").Value = '3.2 °C 4:59 TU'
$ws.Range("O15").Value = '6.9 °C'
$ws.Range("E16").Value = '2026-02-06 05:18:18'
$ws.Range("E17").Value = '2026-02-06 05:18:20'
$ws.Range("J17").Value = '996.8 hPa'
$ws.Range("O17").Value = '3.1 °C'
$ws.Range("E18").Value = '2026-02-06 05:18:23'
$ws.Range("N18").Value = '-5.4 °C 4:38 TU'
$ws.Range("E19").Value = '2026-02-06 05:18:26'
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = '97%'
$ws.Range("I19").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("J19").Value = '997.2 hPa'
$ws.Range("E20").Value = '2026-02-06 05:18:28'
$ws.Range("O20").Value = '-2.1 °C'
$ws.Range("E21").Value = '2026-02-06 05:18:31'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '86%'
$ws.Range("I21").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("J21").Value = '994.5 hPa'
$ws.Range("N21").Value = '2.9 °C 4:56 TU'
$ws.Range("O21").Value = '4.9 °C'
$ws.Range("E22").Value = '2026-02-06 05:18:33'
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = '82%'
$ws.Range("I22").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("K22").Value = '-0.1 MJ/m2'
$ws.Range("N22").Value = '4.3 °C 4:59 TU'
$ws.Range("O22").Value = '8.1 °C'
$ws.Range("E23").Value = '2026-02-06 05:18:35'
$ws.Range("J23").Value = '993.8 hPa'
$ws.Range("E24").Value = '2026-02-06 05:18:38'
$ws.Range("J24").Value = '992.6 hPa'
$ws.Range("K24").Value = '-0.1 MJ/m2'
$ws.Range("E25").Value = '2026-02-06 05:18:40'
$ws.Range("J25").Value = '995.9 hPa'
$ws.Range("N25").Value = '0.8 °C 4:35 TU'
$ws.Range("O25").Value = '2.1 °C'
$ws.Range("E26").Value = '2026-02-06 05:18:42'
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = '81%'
$ws.Range("I26").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("N26").Value = '-1.4 °C 4:45 TU'
$ws.Range("E27").Value = '2026-02-06 05:18:45'
$ws.Range("J27").Value = '993.5 hPa'
$ws.Range("N27").Value = '5.3 °C 4:33 TU'
$ws.Range("O27").Value = '7.4 °C'
$ws.Range("E28").Value = '2026-02-06 05:18:47'
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '89%'
$ws.Range("I28").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("J28").Value = '996.5 hPa'
$ws.Range("N28").Value = '0.2 °C 4:59 TU'
$ws.Range("O28").Value = '2.9 °C'
$ws.Range("E29").Value = '2026-02-06 05:18:50'
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = '63%'
$ws.Range("I29").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("K29").Value = '-0.1 MJ/m2'
$ws.Range("O29").Value = '11.6 °C'
$ws.Range("E30").Value = '2026-02-06 05:18:52'
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '77%'
$ws.Range("I30").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("N30").Value = '-6.4 °C 4:52 TU'
$ws.Range("O30").Value = '-3.7 °C'
$ws.Range("E31").Value = '2026-02-06 05:18:54'
$ws.Range("J31").Value = '996.7 hPa'
$ws.Range("N31").Value = '3.9 °C 4:42 TU'
$ws.Range("O31").Value = '5.0 °C'
$ws.Range("E32").Value = '2026-02-06 05:18:57'
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '51%'
$ws.Range("I32").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("J32").Value = '995.1 hPa'
$ws.Range("N32").Value = '11.5 °C 4:43 TU'
$ws.Range("O32").Value = '14.9 °C'
$ws.Range("E33").Value = '2026-02-06 05:18:59'
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '97%'
$ws.Range("I33").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("N33").Value = '5.0 °C 4:43 TU'
$ws.Range("O33").Value = '6.7 °C'
$ws.Range("E34").Value = '2026-02-06 05:19:01'
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = '78%'
$ws.Range("I34").Copy()
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("N34").Value = '3.1 °C 4:59 TU'
$ws.Range("O34").Value = '7.9 °C'
$ws.Range("E35").Value = '2026-02-06 05:19:04'
$ws.Range("N35").Value = '-3.4 °C 4:51 TU'
$ws.Range("O35").Value = '-3.1 °C'
$ws.Range("E36").Value = '2026-02-06 05:19:06'
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = '67%'
$ws.Range("I36").Copy()
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("J36").Value = '996.6 hPa'
$ws.Range("N36").Value = '9.2 °C 4:58 TU'
$ws.Range("O36").Value = '11.6 °C'
$excel.CutCopyMode = $false
